$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily-push row appended below the existing data (row 38 -> row 39).
# The date column holds plain text like "2025/09/30" (same as every other
# row), not a real date serial, so force text entry with a leading
# apostrophe to stop Excel's automatic date recognition from converting it.
$ws.Cells.Item(39, 1).Value = "'2025/09/30"
$ws.Cells.Item(39, 2).Value = "火"
$ws.Cells.Item(39, 3).Value = 9
$ws.Cells.Item(39, 4).Value = 172
